$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A1").Value = 3.545942068099976
$ws.Range("B1").Value = 2.927317142486572
$ws.Range("C1").Value = 15
$ws.Range("D1").Value = 2.932314157485962
$ws.Range("E1").Value = 3.049896717071533
